# Update countries & provincias Spain
# - Chile overtakes Dinamarca in the ranking (row 28 becomes Chile w/ new
#   totals, row 29 becomes Dinamarca carrying the totals Chile's old row had)
# - Refresh several countries' case counts (Estados Unidos, Alemania,
#   Austria, Pakistan, Kazajistan, Eslovaquia, Birmania)
# - Bump the "Datos actualizados" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28/29 swap: Chile moves above Dinamarca -----------------------
$ws.Range("A28").Value = "Chile"
$ws.Range("B28").Value = 5972
$ws.Range("C28").Value = 426
$ws.Range("D28").Value = 1274
$ws.Range("E28").Value = 4641
$ws.Range("F28").Value = 360
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = 57

$ws.Range("A29").Value = "Dinamarca"
$ws.Range("B29").Value = 5635
$ws.Range("C29").Value = 233
$ws.Range("D29").Value = 1736
$ws.Range("E29").Value = 3662
$ws.Range("F29").Value = 120
$ws.Range("G29").Value = 19
$ws.Range("H29").Value = 237

# --- Row 4: Estados Unidos ----------------------------------------------
$ws.Range("B4").Value = 435780
$ws.Range("C4").Value = 853
$ws.Range("D4").Value = 22941
$ws.Range("E4").Value = 397974
$ws.Range("G4").Value = 77
$ws.Range("H4").Value = 14865

# --- Row 7: Alemania ------------------------------------------------------
$ws.Range("B7").Value = 113615
$ws.Range("C7").Value = 319
$ws.Range("E7").Value = 64966

# --- Row 19: Austria --------------------------------------------------
$ws.Range("B19").Value = 13163
$ws.Range("C19").Value = 221
$ws.Range("E19").Value = 7628

# --- Row 34: Pakistan ---------------------------------------------------
$ws.Range("B34").Value = 4462
$ws.Range("C34").Value = 199
$ws.Range("E34").Value = 3827

# --- Row 75: Kazajistan -------------------------------------------------
$ws.Range("D75").Value = 60
$ws.Range("E75").Value = 697

# --- Row 78: Eslovaquia -------------------------------------------------
$ws.Range("D78").Value = 23
$ws.Range("E78").Value = 676

# --- Row 162: Birmania ----------------------------------------------------
$ws.Range("B162").Value = 23
$ws.Range("C162").Value = 1
$ws.Range("E162").Value = 19

# --- Timestamp banner -----------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 16:22"
